$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the "_GoBack" bookmark from the "Confirm expected log messages"
#    intro paragraph to the very start of the document (start of "TEST 1").
# ---------------------------------------------------------------------------

# Remove the existing _GoBack bookmark (currently sitting between the
# " (testing failure conditions)" run and the ":" run).
$goBack = $d.Bookmarks("_GoBack")
$goBack.Delete()

# Re-create it collapsed at the very beginning of the document. A
# zero-length range exactly at document position 0 cannot be used directly
# to add a bookmark reliably, so nudge it: insert a placeholder character at
# position 0, add the bookmark right after it (a real, non-degenerate
# position), then remove the placeholder again. The bookmark stays anchored
# at the true start of the document, immediately before the "TEST 1" run.
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("X")
$afterPlaceholder = $d.Range(1, 1)
$d.Bookmarks.Add("_GoBack", $afterPlaceholder)
$placeholder = $d.Range(0, 1)
$placeholder.Delete()

# ---------------------------------------------------------------------------
# 2) Consolidate the runs that make up each of the log-message paragraphs in
#    TEST 3 into a single run apiece (the visible text is unchanged; only
#    the run/formatting boundaries collapse together). Using Find & Replace
#    on the whole paragraph text merges the matched runs into one.
# ---------------------------------------------------------------------------

function Merge-ParagraphRuns($text) {
    $null = $d.Content.Find.Execute($text, $true, $false, $false, $false, $false, $true, 1, $false, $text, 2)
}

Merge-ParagraphRuns "WARNING: (TEST_ASSERT_COMPLETE_REFDS) User must ensure PhUSE/CSS utilities are in the AUTOCALL path."
Merge-ParagraphRuns "ERROR: (ASSERT_COMPLETE_REFDS) Result is FAIL. Obs missing from reference dset MY_REFERENCE: 1.5 in_ds1=0 in_ds2=1"
Merge-ParagraphRuns "ERROR: (ASSERT_COMPLETE_REFDS) Result is FAIL. Obs missing from reference dset MY_REFERENCE_C:  Record 1.5 in_ds1=0 in_ds2=1"
Merge-ParagraphRuns "ERROR: (ASSERT_COMPLETE_REFDS) Result is FAIL. Obs missing from reference dset MY_REF_2: 1.002 Rec D Subrec 1.002 in_ds1=0 in_ds2=0 in_ds3=1 in_ds4=0"
Merge-ParagraphRuns "ERROR: (ASSERT_COMPLETE_REFDS) Result is FAIL. Obs missing from reference dset MY_REF_2: 1.5 Record B Subrec 0.001 in_ds1=0 in_ds2=1 in_ds3=0 in_ds4=0"
Merge-ParagraphRuns "ERROR: (ASSERT_COMPLETE_REFDS) Result is FAIL. Obs missing from reference dset MY_REF_2: 2.003 Record C Subrec 400 in_ds1=0 in_ds2=0 in_ds3=0 in_ds4=1"
